$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 with date, version name and change description
$ws.Range("A5").Value2 = 44729
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("B5").Value = "0616YjosetorresMAN11"
$ws.Range("C5").Value = "Completa los calculos de vitamina y corrige % de macronutrientes en dieta de 500 gr"

$ws.Range("C5").Select()
